# The table currently occupies columns B:F (with an extra, stray column A
# holding a style-only helper column). The target layout shifts every
# column one slot to the left (B->A, C->B, D->C, E->D, F->E) and drops the
# old column A entirely. Deleting column A achieves exactly that shift.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").Delete()
